$wb = $excel.ActiveWorkbook

# Hunk: @@ -2240,22 +2240,22 @@
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 1756
$ws.Range("I33").Value = 1870
$ws.Range("K33").Value = 1870
$ws.Range("M33").Value = -1641

# Hunk: @@ -2479,25 +2479,25 @@
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 3075.6155
$ws.Range("I38").Value = 1599.8
$ws.Range("J38").Value = 3998
$ws.Range("K38").Value = 4799.4
$ws.Range("L38").Value = 11994
$ws.Range("M38").Value = -4427.4
$ws.Range("N38").Value = -12738

# Hunk: @@ -5886,22 +5886,22 @@
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 1044.1428
$ws.Range("I106").Value = 1044.1428
$ws.Range("K106").Value = 1044.1428
$ws.Range("M106").Value = -413.1428000000001

# Hunk: @@ -7313,25 +7313,22 @@
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 2001359.8
$ws.Range("I135").Value = 2001359.8
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 18012238.2
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -18009703.2
$ws.Range("N135").ClearContents()

# Hunk: @@ -7414,25 +7411,25 @@
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2379.1333
$ws.Range("I137").Value = 2273.7917
$ws.Range("J137").Value = 2800.5
$ws.Range("K137").Value = 6821.375100000001
$ws.Range("L137").Value = 8401.5
$ws.Range("M137").Value = -4271.375100000001
$ws.Range("N137").Value = -13501.5

# Hunk: @@ -7466,22 +7463,22 @@
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3399.573
$ws.Range("I138").Value = 1000.025
$ws.Range("K138").Value = 3000.075
$ws.Range("M138").Value = 2139.925

# Hunk: @@ -9232,22 +9229,22 @@
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1671237
$ws.Range("I32").Value = 1815529.2
$ws.Range("K32").Value = 1815529.2
$ws.Range("M32").Value = -1815242.2

# Hunk: @@ -15544,25 +15541,25 @@
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 7938173.5
$ws.Range("I20").Value = 20835234
$ws.Range("J20").Value = 1520.6923
$ws.Range("K20").Value = 20835234
$ws.Range("L20").Value = 1520.6923
$ws.Range("M20").Value = -20834987
$ws.Range("N20").Value = -2014.6923

# Hunk: @@ -19376,22 +19373,22 @@
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 3499316.5
$ws.Range("I99").Value = 2606
$ws.Range("K99").Value = 2606
$ws.Range("M99").Value = -1108

# Hunk: @@ -21052,22 +21049,22 @@
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4831
$ws.Range("I134").Value = 1382.9656
$ws.Range("K134").Value = 4148.8968
$ws.Range("M134").Value = -1613.8968

# Hunk: @@ -22947,25 +22944,25 @@
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9573.969999999999
$ws.Range("I31").Value = 4499.727
$ws.Range("J31").Value = 12111.091
$ws.Range("K31").Value = 4499.727
$ws.Range("L31").Value = 12111.091
$ws.Range("M31").Value = -4204.727
$ws.Range("N31").Value = -12701.091

# Hunk: @@ -23094,25 +23091,25 @@
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 9573.969999999999
$ws.Range("I34").Value = 4499.727
$ws.Range("J34").Value = 12111.091
$ws.Range("K34").Value = 4499.727
$ws.Range("L34").Value = 12111.091
$ws.Range("M34").Value = -4297.727
$ws.Range("N34").Value = -12515.091

# Hunk: @@ -24261,22 +24258,22 @@
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 11116544
$ws.Range("I58").Value = 17242932
$ws.Range("K58").Value = 17242932
$ws.Range("M58").Value = -17242729

# Hunk: @@ -26267,25 +26264,25 @@
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 9506.77
$ws.Range("I99").Value = 14447
$ws.Range("J99").Value = 7311.1113
$ws.Range("K99").Value = 14447
$ws.Range("L99").Value = 7311.1113
$ws.Range("M99").Value = -12949
$ws.Range("N99").Value = -10307.1113

# Hunk: @@ -27587,25 +27584,25 @@
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 9506.77
$ws.Range("I126").Value = 14447
$ws.Range("J126").Value = 7311.1113
$ws.Range("K126").Value = 43341
$ws.Range("L126").Value = 21933.3339
$ws.Range("M126").Value = -40871
$ws.Range("N126").Value = -26873.3339

# Hunk: @@ -27878,22 +27875,22 @@
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 5435.048
$ws.Range("I132").Value = 2370.652
$ws.Range("K132").Value = 7111.956
$ws.Range("M132").Value = -4581.956

# Hunk: @@ -27979,25 +27976,25 @@
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 5845.325
$ws.Range("I134").Value = 1388.6471
$ws.Range("J134").Value = 9139.392
$ws.Range("K134").Value = 4165.9413
$ws.Range("L134").Value = 27418.176
$ws.Range("M134").Value = -1630.9413
$ws.Range("N134").Value = -32488.176

# Hunk: @@ -28031,22 +28028,19 @@
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

# Hunk: @@ -28080,22 +28074,22 @@
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 11116544
$ws.Range("I136").Value = 17242932
$ws.Range("K136").Value = 51728796
$ws.Range("M136").Value = -51726246

# Hunk: @@ -28630,25 +28624,25 @@
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 3079530.8
$ws.Range("I5").Value = 5000550
$ws.Range("J5").Value = 5900
$ws.Range("K5").Value = 15001650
$ws.Range("L5").Value = 17700
$ws.Range("M5").Value = -15001538
$ws.Range("N5").Value = -17924

# Hunk: @@ -28985,25 +28979,25 @@
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 129.16667
$ws.Range("I12").Value = 5
$ws.Range("J12").Value = 140.45454
$ws.Range("K12").Value = 15
$ws.Range("L12").Value = 421.36362
$ws.Range("M12").Value = 158
$ws.Range("N12").Value = -767.3636200000001

# Hunk: @@ -29086,22 +29080,22 @@
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 12820853
$ws.Range("I14").Value = 12820853
$ws.Range("K14").Value = 38462559
$ws.Range("M14").Value = -38462386

# Hunk: @@ -31195,22 +31189,22 @@
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 7498.3335
$ws.Range("I56").Value = 7498.3335
$ws.Range("K56").Value = 7498.3335
$ws.Range("M56").Value = -6968.3335

# Hunk: @@ -33322,25 +33316,25 @@
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 3855.9
$ws.Range("J98").Value = 4222.857
$ws.Range("L98").Value = 12668.571
$ws.Range("N98").Value = -15664.571

# Hunk: @@ -33760,25 +33754,25 @@
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 16667018
$ws.Range("J107").Value = 22222558
$ws.Range("L107").Value = 66667674
$ws.Range("N107").Value = -66671514

# Hunk: @@ -34868,25 +34862,25 @@
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 142858720
$ws.Range("I129").Value = 965
$ws.Range("J129").Value = 200001820
$ws.Range("K129").Value = 2895
$ws.Range("L129").Value = 600005460
$ws.Range("M129").Value = 2105
$ws.Range("N129").Value = -600015460

# Hunk: @@ -34969,25 +34963,25 @@
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1346
$ws.Range("J131").Value = 1430.619
$ws.Range("L131").Value = 4291.857
$ws.Range("N131").Value = -14371.857

# Hunk: @@ -35174,25 +35168,25 @@
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 3079530.8
$ws.Range("I135").Value = 5000550
$ws.Range("J135").Value = 5900
$ws.Range("K135").Value = 45004950
$ws.Range("L135").Value = 53100
$ws.Range("M135").Value = -45002415
$ws.Range("N135").Value = -58170

# Hunk: @@ -40206,25 +40200,25 @@
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1011
$ws.Range("I97").Value = 922.59375
$ws.Range("J97").Value = 1325.3334
$ws.Range("K97").Value = 922.59375
$ws.Range("L97").Value = 1325.3334
$ws.Range("M97").Value = -426.59375
$ws.Range("N97").Value = -2317.3334

# Hunk: @@ -41915,25 +41909,25 @@
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4229.8057
$ws.Range("I132").Value = 1448.76
$ws.Range("J132").Value = 10550.363
$ws.Range("K132").Value = 4346.28
$ws.Range("L132").Value = 31651.089
$ws.Range("M132").Value = -1816.28
$ws.Range("N132").Value = -36711.089

# Hunk: @@ -43182,22 +43176,22 @@
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 870.9167
$ws.Range("I16").Value = 759.1818
$ws.Range("K16").Value = 759.1818
$ws.Range("M16").Value = -589.1818

# Hunk: @@ -43479,25 +43473,25 @@
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1700.3667
$ws.Range("I22").Value = 1259.2593
$ws.Range("J22").Value = 5670.3335
$ws.Range("K22").Value = 1259.2593
$ws.Range("L22").Value = 5670.3335
$ws.Range("M22").Value = -964.2592999999999
$ws.Range("N22").Value = -6260.3335

# Hunk: @@ -43718,25 +43712,25 @@
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 1700.3667
$ws.Range("I27").Value = 1259.2593
$ws.Range("J27").Value = 5670.3335
$ws.Range("K27").Value = 1259.2593
$ws.Range("L27").Value = 5670.3335
$ws.Range("M27").Value = -1152.2593
$ws.Range("N27").Value = -5884.3335

# Hunk: @@ -48770,22 +48764,22 @@
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 12507225
$ws.Range("I132").Value = 26318842
$ws.Range("K132").Value = 78956526
$ws.Range("M132").Value = -78953996

# Hunk: @@ -48966,25 +48960,25 @@
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 11931.682
$ws.Range("I136").Value = 2874.5
$ws.Range("J136").Value = 13944.389
$ws.Range("K136").Value = 8623.5
$ws.Range("L136").Value = 41833.167
$ws.Range("M136").Value = -6073.5
$ws.Range("N136").Value = -46933.167

# Hunk: @@ -53132,22 +53126,22 @@
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 28858798
$ws.Range("I81").Value = 401318.2
$ws.Range("K81").Value = 802636.4
$ws.Range("M81").Value = -801575.4

# Hunk: @@ -53276,22 +53270,22 @@
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 28858798
$ws.Range("I84").Value = 401318.2
$ws.Range("K84").Value = 4013182
$ws.Range("M84").Value = -4007878

# Hunk: @@ -54400,25 +54394,25 @@
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 13889618
$ws.Range("J107").Value = 47620244
$ws.Range("L107").Value = 142860732
$ws.Range("N107").Value = -142864572

# Hunk: @@ -54682,22 +54676,22 @@
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 802.8946999999999
$ws.Range("I113").Value = 676
$ws.Range("K113").Value = 2028
$ws.Range("M113").Value = 142

# Hunk: @@ -55619,22 +55613,22 @@
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 25020794
$ws.Range("I132").Value = 45465404
$ws.Range("K132").Value = 136396212
$ws.Range("M132").Value = -136393682

# Hunk: @@ -55815,25 +55809,25 @@
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 31287318
$ws.Range("I136").Value = 83334620
$ws.Range("J136").Value = 58933.65
$ws.Range("K136").Value = 250003860
$ws.Range("L136").Value = 176800.95
$ws.Range("M136").Value = -250001310
$ws.Range("N136").Value = -181900.95
